$wb = $excel.ActiveWorkbook

# --- Text update: "Ready for handoff" -> "In Translation" -------------
# Overview sheet tracks per-language status in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn / de-de detail sheets carry the same status in column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width update (status columns got narrower to fit new text) -
# (ColumnWidth snaps to the nearest 1/6-character pixel grid; 12.58 is the
# input that lands on the grid point closest to the target 13.41 width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.58
$wsOverview.Columns.Item(6).ColumnWidth = 12.58

$wsZhCn.Columns.Item(3).ColumnWidth = 12.58

$wsDeDe.Columns.Item(3).ColumnWidth = 12.58
